$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Agosto de 2020 a las 14:38"

# Update country names where ranking order shifted due to refreshed data
$ws.Cells.Item(63, 1).Value = "Uzbekistan"
$ws.Cells.Item(64, 1).Value = "Moldavia"
$ws.Cells.Item(79, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(80, 1).Value = "Estado de Palestina"
$ws.Cells.Item(143, 1).Value = "Uganda"
$ws.Cells.Item(144, 1).Value = "Georgia"
$ws.Cells.Item(157, 1).Value = "Lesoto"
$ws.Cells.Item(158, 1).Value = "Crucero"

# Update numeric statistics (columns B:H) for rows affected by the data refresh
# Row 4
$ws.Cells.Item(4, 2).Value = 4765155
$ws.Cells.Item(4, 3).Value = 837
$ws.Cells.Item(4, 4).Value = 2363165
$ws.Cells.Item(4, 5).Value = 2244069
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 23
$ws.Cells.Item(4, 8).Value = 157921

# Row 24
$ws.Cells.Item(24, 2).Value = 129151
$ws.Cells.Item(24, 3).Value = 2447
$ws.Cells.Item(24, 4).Value = 91886
$ws.Cells.Item(24, 5).Value = 32397
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 63
$ws.Cells.Item(24, 8).Value = 4868

# Row 40
$ws.Cells.Item(40, 2).Value = 68067
$ws.Cells.Item(40, 3).Value = 121
$ws.Cells.Item(40, 4).Value = 62896
$ws.Cells.Item(40, 5).Value = 4604
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 4
$ws.Cells.Item(40, 8).Value = 567

# Row 47
$ws.Cells.Item(47, 2).Value = 51463
$ws.Cells.Item(47, 3).Value = 153
$ws.Cells.Item(47, 4).Value = 36984
$ws.Cells.Item(47, 5).Value = 12741
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 1
$ws.Cells.Item(47, 8).Value = 1738

# Row 63
$ws.Cells.Item(63, 2).Value = 25336
$ws.Cells.Item(63, 3).Value = 553
$ws.Cells.Item(63, 4).Value = 15833
$ws.Cells.Item(63, 5).Value = 9352
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 4
$ws.Cells.Item(63, 8).Value = 151

# Row 64
$ws.Cells.Item(64, 2).Value = 25113
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(64, 4).Value = 17816
$ws.Cells.Item(64, 5).Value = 6508
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 789

# Row 79
$ws.Cells.Item(79, 2).Value = 12296
$ws.Cells.Item(79, 3).Value = 420
$ws.Cells.Item(79, 4).Value = 6312
$ws.Cells.Item(79, 5).Value = 5632
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 13
$ws.Cells.Item(79, 8).Value = 352

# Row 80
$ws.Cells.Item(80, 2).Value = 12160
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(80, 4).Value = 5324
$ws.Cells.Item(80, 5).Value = 6753
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 83

# Row 84
$ws.Cells.Item(84, 2).Value = 11054
$ws.Cells.Item(84, 3).Value = 163
$ws.Cells.Item(84, 4).Value = 6882
$ws.Cells.Item(84, 5).Value = 3675
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 4
$ws.Cells.Item(84, 8).Value = 497

# Row 100
$ws.Cells.Item(100, 2).Value = 5260
$ws.Cells.Item(100, 3).Value = 36
$ws.Cells.Item(100, 4).Value = 4373
$ws.Cells.Item(100, 5).Value = 738
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 4
$ws.Cells.Item(100, 8).Value = 149

# Row 118
$ws.Cells.Item(118, 2).Value = 2816
$ws.Cells.Item(118, 3).Value = 1
$ws.Cells.Item(118, 4).Value = 2514
$ws.Cells.Item(118, 5).Value = 291
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 11

# Row 132
$ws.Cells.Item(132, 2).Value = 1907
$ws.Cells.Item(132, 3).Value = 14
$ws.Cells.Item(132, 4).Value = 1825
$ws.Cells.Item(132, 5).Value = 72
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 10

# Row 143
$ws.Cells.Item(143, 2).Value = 1182
$ws.Cells.Item(143, 3).Value = 6
$ws.Cells.Item(143, 4).Value = 1045
$ws.Cells.Item(143, 5).Value = 133
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 4

# Row 144
$ws.Cells.Item(144, 2).Value = 1177
$ws.Cells.Item(144, 3).Value = 6
$ws.Cells.Item(144, 4).Value = 955
$ws.Cells.Item(144, 5).Value = 205
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 17

# Row 157
$ws.Cells.Item(157, 2).Value = 718
$ws.Cells.Item(157, 3).Value = 16
$ws.Cells.Item(157, 4).Value = 173
$ws.Cells.Item(157, 5).Value = 526
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 5
$ws.Cells.Item(157, 8).Value = 19

# Row 158
$ws.Cells.Item(158, 2).Value = 712
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 651
$ws.Cells.Item(158, 5).Value = 48
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 13

# Row 195
$ws.Cells.Item(195, 2).Value = 55
$ws.Cells.Item(195, 3).Value = 1
$ws.Cells.Item(195, 4).Value = 45
$ws.Cells.Item(195, 5).Value = 10
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

